# Daily attendance processing - 2025-12-08 08:37:38
# Re-orders the "Recorded By" attendance-taker lists for several sessions
# (System/clock-in order changed upstream) and refreshes the derived
# Missing/Pending session counters, plus flips session PARASITOLOGY #6
# (row 18) from "Pending" to "Not Recorded" now that its scheduled time
# has passed without being recorded.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder "Recorded By" (column G) email lists ---------------------
$ws.Range("G2").Value  = "System, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G3").Value  = "hend_mahmoud@med.asu.edu.eg, System, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G4").Value  = "hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, servinaz@med.asu.edu.eg"
$ws.Range("G5").Value  = "eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G7").Value  = "Kerelos.zareef@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg"
$ws.Range("G9").Value  = "Safa.hany@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"
$ws.Range("G12").Value = "Marina.youhana@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, dina.adel@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg"
$ws.Range("G15").Value = "mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"
$ws.Range("G24").Value = "youstina.gamil@med.asu.edu.eg, Sarah.Mahdy@med.asu.edu.eg"
$ws.Range("G27").Value = "hana.amr@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"
$ws.Range("G28").Value = "Aya_hamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"
$ws.Range("G30").Value = "aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"

# --- Refresh derived statistics counters -------------------------------
$ws.Range("L7").Value = 3   # BIOCHEMISTRY LAB/CBL - Missing Sessions
$ws.Range("L8").Value = 3   # BIOCHEMISTRY LAB/CBL - Pending Sessions

$ws.Range("P15").Value = 3  # PARASITOLOGY group row - Missing
$ws.Range("Q15").Value = 3  # PARASITOLOGY group row - Pending

# --- Session PARASITOLOGY #6 (row 18) is now overdue -> Not Recorded ---
# Copy the "Not Recorded" look (pink fill/font) from row 11, which already
# carries that formatting, onto row 18 without disturbing row 18's own
# Year/Group/Subject/Session/Date/Time/Students values.
$ws.Range("A11:I11").Copy() | Out-Null
$ws.Range("A18:I18").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("I18").Value = "Not Recorded"
